$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-15 Wednesday" "2024-05-16 Thursday"
Replace-Text "32×14=" "73×60="
Replace-Text "24×90=" "73×35="
Replace-Text "19×41=" "34×71="
Replace-Text "42×47=" "13×27="
Replace-Text "80×38=" "31×22="
Replace-Text "66×70=" "71×41="
Replace-Text "76×18=" "96×42="
Replace-Text "25×75=" "20×67="
Replace-Text "66×95=" "56×48="
Replace-Text "82×61=" "25×99="
Replace-Text "29×87=" "97×12="
Replace-Text "57×94=" "30×65="
Replace-Text "75×75=" "60×11="
Replace-Text "20×28=" "33×98="
Replace-Text "65×81=" "18×68="
Replace-Text "18×13=" "87×45="
Replace-Text "27×62=" "16×27="
Replace-Text "97×49=" "83×24="
Replace-Text "86×67=" "39×34="
Replace-Text "96×29=" "97×90="
Replace-Text "33×11=" "65×78="
Replace-Text "55×66=" "30×18="
Replace-Text "87×28=" "33×80="
Replace-Text "97×19=" "35×20="
Replace-Text "49×85=" "94×93="
